# "Generate Report for Handoff"
#
# The localization-status report tracks one row per source file, per
# language sheet (zh-cn / de-de) plus a roll-up on the Overview sheet.
# The c9e48bb5-... file has just been handed off for (machine) translation,
# so its status / priority / timestamps move forward:
#
#   Status   : "In Translation"        -> "Ready for handoff"
#   Priority : "ht" (human translate)  -> "mt" (machine translate)
#   Latest Handoff Datetime moves forward a little over half a minute.
#
# The Overview sheet mirrors the per-language Status + the newer of the two
# per-language handoff timestamps in its "Latest HO Xliff Generate Date"
# column.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newPriority = "mt"

$zhDatetime = "2016-08-17 18:13:29"
$deDatetime = "2016-08-17 18:13:34"

# --- zh-cn sheet: row 3 is the c9e48bb5 file -------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("E3").Value = $newPriority
$wsZh.Range("H3").Value = $zhDatetime

# --- de-de sheet: row 3 is the c9e48bb5 file -------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("E3").Value = $newPriority
$wsDe.Range("H3").Value = $deDatetime

# --- Overview sheet: row 3 is the c9e48bb5 file ----------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $deDatetime

# --- Widen the "Status" columns so the longer "Ready for handoff" text
#     fits without truncation (Overview!E:F, zh-cn!C, de-de!C).
$newColumnWidth = 98 / 6

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZh.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDe.Columns.Item(3).ColumnWidth = $newColumnWidth
